# "change ui of department and designation + formats"
#
# The sheet originally had:
#   Row 1: "Your Organization Name" (bold label, A1) with an empty B1
#   Row 2: "Date" (bold label, A2) with an empty, date-formatted B2
#   Row 3: the real header row -> Code | Name | Status | Leaving Date | Reason
#
# The edit removes the two throwaway label rows at the top so the header
# row becomes row 1, and moves the active selection to F2 (just past the
# new header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Your Organization Name" / "Date" rows entirely; this shifts
# the header row (Code/Name/Status/Leaving Date/Reason) up from row 3 to
# row 1 and shrinks the used range to A1:E1.
$ws.Rows("1:2").Delete()

# Match the new selection recorded in the saved workbook.
$ws.Range("F2").Select()
